$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.410.59'
$ws.Range('E2').Value = '  -1.42%  '

$ws.Range('D3').Value = '2.371.90'
$ws.Range('E3').Value = '  +4.59%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.77%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.650'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.03%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.63'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +12.18%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.470'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.32%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0979'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.53%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.84'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.79%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '27.10'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.28%  '

$ws.Range('D13').Value = '2.726.89'
$ws.Range('E13').Value = '  +4.62%  '

$ws.Range('E14').Value = '  -0.05%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.98'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.42%  '

$ws.Range('E16').Value = '  +1.71%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.851'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.00%  '

$ws.Range('D18').Value = '2.374.72'
$ws.Range('E18').Value = '  +4.78%  '

$ws.Range('D19').Value = '43.421.77'
$ws.Range('E19').Value = '  -1.24%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000100'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.35%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '74.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.71%  '

$ws.Range('E22').Value = '  +3.28%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.72%  '

$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.08%  '

$ws.Range('B25').Value = 'WEMIXToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.76'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +13.06%  '

$ws.Range('E26').Value = '  +1.13%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.89%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '22.97'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.22%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.01'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.13%  '

$ws.Range('E30').Value = '  +0.30%  '

$ws.Range('E31').Value = '  +4.07%  '

$ws.Range('E32').Value = '  -6.49%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.127'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.28%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.81%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0690'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.55%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.82%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.46'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.23%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.57'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.38%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.66'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.44%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0256'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.08%  '

$ws.Range('B41').Value = 'BinanceUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.05%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.92'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.74%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '18.56'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.60%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.06%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '100.08'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.32%  '

$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.52'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.73%  '

$ws.Range('E47').Value = '  +1.68%  '

$ws.Range('E48').Value = '  +0.29%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000208'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.93%  '

$ws.Range('D50').Value = '1.444.02'
$ws.Range('E50').Value = '  -0.59%  '

$ws.Range('D51').Value = '2.600.18'
$ws.Range('E51').Value = '  +4.87%  '
